$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.439070463180542
$ws.Range("B1").Value = 2.784838914871216
$ws.Range("C1").Value = 1.894824028015137
$ws.Range("D1").Value = 1.645299196243286
$ws.Range("E1").Value = 1.611118674278259
